# Update NATMI ligand-receptor edge table: sending cluster re-run with new
# TPM values changed the "Sending cluster" from ECs to MuSCs (ligand/receptor
# pair Wnt8a/Fzd8 unchanged), and refreshed all derived expression metrics
# for the three target clusters (ECs, FAPs, MuSCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: MuSCs -> ECs -------------------------------------------------
$ws.Cells.Item(2, 1).Value = "MuSCs"   # A2 Sending cluster
$ws.Cells.Item(2, 2).Value = "Wnt8a"   # B2 Ligand symbol
$ws.Cells.Item(2, 3).Value = "Fzd8"    # C2 Receptor symbol
$ws.Cells.Item(2, 4).Value = "ECs"     # D2 Target cluster
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.01490866666666667
$ws.Cells.Item(2, 8).Value = 0.044726
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.083576666666666
$ws.Cells.Item(2, 14).Value = 9.250729999999999
$ws.Cells.Item(2, 15).Value = 0.2272509363535097
$ws.Cells.Item(2, 16).Value = 0.2272509363535097
$ws.Cells.Item(2, 17).Value = 0.04597201666444444
$ws.Cells.Item(2, 18).Value = 0.41374814998
$ws.Cells.Item(2, 19).Value = 0.2272509363535097
$ws.Cells.Item(2, 20).Value = 0.2272509363535097

# --- Row 3: MuSCs -> FAPs -------------------------------------------------
$ws.Cells.Item(3, 1).Value = "MuSCs"   # A3 Sending cluster
$ws.Cells.Item(3, 2).Value = "Wnt8a"   # B3 Ligand symbol
$ws.Cells.Item(3, 3).Value = "Fzd8"    # C3 Receptor symbol
$ws.Cells.Item(3, 4).Value = "FAPs"    # D3 Target cluster
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.01490866666666667
$ws.Cells.Item(3, 8).Value = 0.044726
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 6.453984666666667
$ws.Cells.Item(3, 14).Value = 19.361954
$ws.Cells.Item(3, 15).Value = 0.4756405360586227
$ws.Cells.Item(3, 16).Value = 0.4756405360586227
$ws.Cells.Item(3, 17).Value = 0.09622030606711111
$ws.Cells.Item(3, 18).Value = 0.8659827546040001
$ws.Cells.Item(3, 19).Value = 0.4756405360586227
$ws.Cells.Item(3, 20).Value = 0.4756405360586227

# --- Row 4: MuSCs -> MuSCs -------------------------------------------------
$ws.Cells.Item(4, 1).Value = "MuSCs"   # A4 Sending cluster
$ws.Cells.Item(4, 2).Value = "Wnt8a"   # B4 Ligand symbol
$ws.Cells.Item(4, 3).Value = "Fzd8"    # C4 Receptor symbol
$ws.Cells.Item(4, 4).Value = "MuSCs"   # D4 Target cluster
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.01490866666666667
$ws.Cells.Item(4, 8).Value = 0.044726
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.031477000000001
$ws.Cells.Item(4, 14).Value = 12.094431
$ws.Cells.Item(4, 15).Value = 0.2971085275878677
$ws.Cells.Item(4, 16).Value = 0.2971085275878677
$ws.Cells.Item(4, 17).Value = 0.06010394676733335
$ws.Cells.Item(4, 18).Value = 0.5409355209060001
$ws.Cells.Item(4, 19).Value = 0.2971085275878677
$ws.Cells.Item(4, 20).Value = 0.2971085275878677
